# Connectors.xlsx edit script
# - Switch hotend wiring (rows 9-12, AWG column C) from 22 awg to 20 awg
# - A/B steppers (row 42) now have pigtailed connectors -> add Len (cm) value
# - Add a new "BED-1" row (row 59) for a phantom tee nut / XT 60 connector
# - Update the active selection on the "Fabricated Cables" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fabricated Cables")

# Hotend wiring switched to 20 AWG (column C, "AWG")
$ws.Range("C9").Value = 20
$ws.Range("C10").Value = 20
$ws.Range("C11").Value = 20
$ws.Range("C12").Value = 20

# Fix #42 - A/B steppers now have pigtailed connectors (add a cable length)
$ws.Range("E42").Value = 100

# New row 59 - additional connector entry
$ws.Range("A59").Value = "BED-1"
$ws.Range("B59").Value = "BED"
$ws.Range("C59").Value = 12
$ws.Range("D59").Value = 2
$ws.Range("E59").Value = 35
$ws.Range("F59").Value = "XT 60"
$ws.Range("G59").Value = "Ring"

# Grow the table ("Table3") to include the newly-added row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H59"))

# Update the sheet's view / selection
$ws.Activate()
$ws.Range("C13").Select()
